# Insert two new rows at 448-449 (existing rows 448:539 shift down to 450:541)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A448:A449").EntireRow.Insert()

# ---- New row 448 ----
$ws.Range("A448").Value = 6
$ws.Range("B448").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C448").Value = "Metropolitana"
$ws.Range("D448").Value = 44951
$ws.Range("E448").Value = 13
$ws.Range("F448").Value = "Fruta"
$ws.Range("G448").Value = 100101
$ws.Range("H448").Value = "Berries"
$ws.Range("I448").Value = 100101001
$ws.Range("J448").Value = "Arándano (blue)"
$ws.Range("K448").Value = "Sin especificar"
$ws.Range("L448").Value = "Especial"
$ws.Range("M448").Value = 250
$ws.Range("N448").Value = 2800
$ws.Range("O448").Value = 2800
$ws.Range("P448").Value = 2800
$ws.Range("Q448").Value = "$/bandeja 2 kilos"
$ws.Range("R448").Value = "Provincia de Curicó"
$ws.Range("S448").Value = 1400
$ws.Range("T448").Value = 2

# ---- New row 449 ----
$ws.Range("A449").Value = 6
$ws.Range("B449").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C449").Value = "Metropolitana"
$ws.Range("D449").Value = 44951
$ws.Range("E449").Value = 13
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100101
$ws.Range("H449").Value = "Berries"
$ws.Range("I449").Value = 100101001
$ws.Range("J449").Value = "Arándano (blue)"
$ws.Range("K449").Value = "Sin especificar"
$ws.Range("L449").Value = "Especial"
$ws.Range("M449").Value = 500
$ws.Range("N449").Value = 3000
$ws.Range("O449").Value = 3000
$ws.Range("P449").Value = 3000
$ws.Range("Q449").Value = "$/bandeja 2 kilos"
$ws.Range("R449").Value = "Región del Maule"
$ws.Range("S449").Value = 1500
$ws.Range("T449").Value = 2
